$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Slightly re-tuned column widths (cosmetic re-save drift)
$ws.Columns.Item(3).ColumnWidth = 18.830729166666668
$ws.Columns.Item(4).ColumnWidth = 14.830729166666666
$ws.Columns.Item(5).ColumnWidth = 32.498697916666664
$ws.Columns.Item(6).ColumnWidth = 24.608072916666668
$ws.Columns.Item(7).ColumnWidth = 22.830729166666668
$ws.Columns.Item(8).ColumnWidth = 22.721354166666668

# 1. Task 1.6 description: drop the trailing clause about inventory flow reporting
$ws.Range("E12").Value = "Xây dựng module thống kê hàng hóa trong kho"

# 2. Task 1.9 row: was the "help docs" task, now becomes the "integrate modules" task
#    (content that used to live two rows further down, at task 1.11)
$ws.Range("E15").Value = "Ghép nối các module để nhận được phần mềm hoàn chỉnh"
$ws.Range("G15").Value = "Dũng"

# 3. Task 1.10 row: was the "install docs" task, now becomes the "run tests" task
#    (content that used to live two rows further down, at task 1.12)
$ws.Range("E16").Value = "Chạy kiểm thử phần mềm kết quả"
$ws.Range("G16").Value = "Lào - Campuchia"

# 4/5. The old rows 17 and 18 (tasks 1.11 / 1.12) no longer exist as separate
#      rows - their content moved up into rows 15/16 above, so clear the
#      leftover task number / description / owner cells and drop their border.
$tail = $ws.Range("C17:H18")
$tail.ClearContents()
$tail.Borders.LineStyle = -4142

# Cursor moved on from the old H14 task row to the now-blank D21
$ws.Range("D21").Select() | Out-Null

$wb.Save()
